# Rename the header row (row 1) labels to their underscore-separated
# equivalents, matching the "database and tables created" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "CUT_SPEED"
$ws.Range("C1").Value = "PIERCE_TIME"
$ws.Range("E1").Value = "MATERIAL_COST"
$ws.Range("F1").Value = "COST_SQUARE_FOOT"
$ws.Range("A1").Value = "STAINLESS_STEEL"
$ws.Range("D1").Value = "WEIGHT"
